$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Benoit Danglades' task is refined from "BDD+Back" to "Back" as the project kicks off
$ws.Range("C3").Value = "Back"

# Reflect the cell the author was last working in
$ws.Range("C3").Select()
